# Aug 21, 2015 Update
# Adds new Python / R cheat-sheet rows, tweaks a couple of existing
# cells, and switches the active sheet/selection back to "Python".
#
# NOTE: the Value assignments below are intentionally ordered to match
# the order new shared strings were first introduced, so the rebuilt
# sharedStrings table lines up with the target workbook.

$wb = $excel.ActiveWorkbook

$wsPython = $wb.Worksheets.Item("Python")
$wsR      = $wb.Worksheets.Item("R")

# Row 14 (Python) - new "NOT in" example in columns E/F
$wsPython.Range("E14").Value = "df[~df['col'].isin(['ye','boi'])]"
$wsPython.Range("F14").Value = "#3 is NOT in"

# Row 24 (Python) - category corrected from "Basic" to "Import / Export", plus a note
$wsPython.Range("A24").Value = "Import / Export"
$wsPython.Range("F24").Value = "Add. Args: index = False"

# Row 68 (R) - new "String Character Length" entry
$wsR.Range("A68").Value = "Basic"
$wsR.Range("B68").Value = "String Character Length"
$wsR.Range("C68").Value = "transform(df, Length=nchar(as.character(A)))"

# Row 69 (R) - new "Approximate Match" entry
$wsR.Range("B69").Value = "Approximate Match"
$wsR.Range("C69").Value = "agrep('string',df`$A)"
$wsR.Range("F69").Value = "Returns location"

# Row 70 (Python) - new "Count Uniques within Dimension" entry
$wsPython.Range("B70").Value = "Count Uniques within Dimension"
$wsPython.Range("C70").Value = "df.groupby('A').CountItem.nunique()"

# Row 71 (Python) - new "Conver to Pandas DF" entry
$wsPython.Range("B71").Value = "Conver to Pandas DF"

# Row 70 (R) - new "Count IF" entry
$wsR.Range("B70").Value = "Count IF"
$wsR.Range("C70").Value = "sum(df`$a=='ye',na.rm=T)"

# Row 72 (Python) - new "Categorize Based on Value" entry
$wsPython.Range("B72").Value = "Categorize Based on Value"
$wsPython.Range("C72").Value = "df['b'] = np.where(df['a'] > 5,'high','low')"

# Row 71 (R) - new "Match items on list" entry
$wsR.Range("B71").Value = "Match items on list"
$wsR.Range("C71").Value = "match(df`$a,matchtable)"
$wsR.Range("F71").Value = "nomatch = 0"

# Row 72 (R) - new "Change values conditionally" entry
$wsR.Range("C72").Value = "df`$a[which(df`$a>0)]<-1"
$wsR.Range("B72").Value = "Change values conditionally"

# Row 73 (Python) - new "Match (Intersection)" entry
$wsPython.Range("B73").Value = "Match (Intersection)"
$wsPython.Range("C73").Value = "set(df1) & set(df2)"

# Row 49 (Python) - additional note about the regex argument
$wsPython.Range("F49").Value = "Add. Args: regex = True (for part of a string)"

# Row 74 (Python) - new "Keep Columns" entry
$wsPython.Range("B74").Value = "Keep Columns"
$wsPython.Range("C74").Value = "df[['a','b']]"

# Row 75 (Python) - new "Keep first occurrence" entry
$wsPython.Range("B75").Value = "Keep first occurrence"
$wsPython.Range("C75").Value = "df.groupby('a').first()"

# ---------------------------------------------------------------------
# View state - "R" keeps its own last selection, but "Python" becomes
# the active/selected sheet again (select R's cell first, then finish
# on Python so Python ends up as the active tab).
# ---------------------------------------------------------------------

$wsR.Range("A68").Select() | Out-Null
$wsPython.Range("C26").Select() | Out-Null
